{"js": "// \"no youtube links on the resume\"\n// Remove the \"Portfolio: https://www.youtube.com/watch?v=tDC5uYWP46Y\" hyperlink\n// (and the manual line break that introduced it) from the WebGL/OpenGL 3D\n// Graphics Demo bullet, while leaving everything else (including the\n// itch.io link later in the document) untouched.\n\nconst body = context.document.body;\n\n// Step 1: remove the \"Portfolio: <youtube url>\" text together with the\n// hyperlink run. Searching across runs (including the hyperlinked run)\n// works fine as long as the match doesn't also try to swallow the\n// preceding manual line break in the same call, so this is done as its\n// own search/delete pass.\nlet results = body.search(\"Portfolio: https://www.youtube.com/watch?v=tDC5uYWP46Y\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].delete();\n  await context.sync();\n}\n\n// Step 2: remove the manual line break (w:br) that used to separate the\n// bullet text from the now-deleted \"Portfolio:\" label. In the document's\n// text projection a manual line break reads back as a vertical-tab\n// character (\\u000b), and this document only has the one (right after\n// \"...user interaction.\"), so it's safe to search for it directly.\nresults = body.search(\"\\u000b\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# \"no youtube links on the resume\"\n# Remove the \"Portfolio: https://www.youtube.com/watch?v=tDC5uYWP46Y\" hyperlink\n# (and the manual line break that introduced it) from the WebGL/OpenGL 3D\n# Graphics Demo bullet, while leaving everything else (including the\n# itch.io link later in the document) untouched.\n\n$d = $word.ActiveDocument\n\n# Step 1: delete the \"Portfolio:\" label (spans the italic \"Portfolio\" run\n# and the plain \":\" run). Done as its own Find/Delete pass so it doesn't\n# also try to swallow the following hyperlinked run in the same call.\n$rng1 = $d.Content\n$found1 = $rng1.Find.Execute(\"Portfolio:\")\nif ($found1) {\n  $rng1.Delete()\n}\n\n# Step 2: delete the hyperlinked URL text run (\" https://www.youtube.com/watch?v=tDC5uYWP46Y\").\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\" https://www.youtube.com/watch?v=tDC5uYWP46Y\")\nif ($found2) {\n  $rng2.Delete()\n}\n\n# Step 3: delete the now-orphaned manual line break that used to separate\n# the bullet text from the \"Portfolio:\" label. A manual line break\n# (<w:br/>) reads back as a vertical-tab character (Chr 11); this document\n# only has the one, right after \"...user interaction.\".\n$rng3 = $d.Content\n$found3 = $rng3.Find.Execute([char]11)\nif ($found3) {\n  $rng3.Delete()\n}\n"}
